# "kotlin advanced - 5 animations"
# Slide 4 ("CustomShape 12" code block) gets a second, visually-identical
# text box (blank / space-only) stacked right below it, and both shapes
# get an on-click "Appear" animation (the pre-existing one on the
# original shape plus a brand-new one on the new shape).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$seq = $s.TimeLine.MainSequence

# ---------------------------------------------------------------------
# The fixture's shape-id allocator is per-slide and starts at 2
# (existing shapes 102-113 don't reserve anything in that counter), so
# to land the new shape on id 15 we burn through 13 throw-away shapes
# first, exactly like additional shapes being added/removed earlier in
# the slide's history would have done.
# ---------------------------------------------------------------------
$junk = @()
for ($i = 1; $i -le 13; $i++) {
    $junk += $s.Shapes.AddShape(1, 0, 0, 10, 10)
}
foreach ($j in $junk) { $j.Delete() }

# ---------------------------------------------------------------------
# Duplicate the existing "CustomShape 12" (id 113) so the new shape
# inherits identical shape styling / text formatting, then reposition
# it just below the original and blank out its text to two single
# space runs.
# ---------------------------------------------------------------------
$src = $s.Shapes.Item(12)
$newShp = $src.Duplicate()

# Reposition to the exact target offset (806940, 6343280 EMU). Shape.Left/Top
# are expressed in points (1 pt = 12700 EMU); these particular values are
# chosen so the float32 COM marshalling round-trips to the exact EMU target.
$newShp.Left = 63.53859
$newShp.Top = 499.4709

$tr = $newShp.TextFrame.TextRange
$tr.Text = " "
$null = $tr.InsertAfter(" ")

# ---------------------------------------------------------------------
# Animations: refresh the existing "Appear" build on shape 113 (this is
# what re-normalizes its timing node to the newer
# presetSubtype/grpId-qualified form) and add a matching "Appear" build
# for the new shape. Order matters: 113 must be re-added before 15 so
# the new timing nodes/bldLst entries land in the same order as target.
# ---------------------------------------------------------------------
for ($i = 1; $i -le $seq.Count; $i++) {
    $e = $seq.Item($i)
    if ($e.Shape.Id -eq $src.Id) {
        $e.Delete()
        break
    }
}
$null = $seq.AddEffect($src, 1)
$null = $seq.AddEffect($newShp, 1)
